$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking values (e.g. "1.004", "17.50")
# are not reinterpreted/rounded as numbers when the new values are assigned.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.907.12'
$ws.Range("E2").Value = '  +1.44%  '

$ws.Range("D3").Value = '1.710.86'
$ws.Range("E3").Value = '  +1.46%  '

$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.53%  '

$ws.Range("D5").Value = '315.51'
$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.62%  '

$ws.Range("D7").Value = '0.4030'
$ws.Range("E7").Value = '  +3.46%  '

$ws.Range("D8").Value = '0.4054'
$ws.Range("E8").Value = '  +0.74%  '

$ws.Range("D9").Value = '1.003'
$ws.Range("E9").Value = '  -0.60%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '1.475'
$ws.Range("E10").Value = '  -0.51%  '

$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = '53.74'
$ws.Range("E11").Value = '  +1.37%  '

$ws.Range("D12").Value = '0.08825'
$ws.Range("E12").Value = '  +1.62%  '

$ws.Range("D13").Value = '26.18'
$ws.Range("E13").Value = '  +7.40%  '

$ws.Range("D14").Value = '7.526'
$ws.Range("E14").Value = '  -0.72%  '

$ws.Range("D15").Value = '8.022'
$ws.Range("E15").Value = '  +1.19%  '

$ws.Range("D16").Value = '0.00001347'
$ws.Range("E16").Value = '  +1.07%  '

$ws.Range("D17").Value = '1.667.52'
$ws.Range("E17").Value = '  -1.19%  '

$ws.Range("D18").Value = '95.52'
$ws.Range("E18").Value = '  -2.73%  '

$ws.Range("D19").Value = '0.07183'
$ws.Range("E19").Value = '  +1.17%  '

$ws.Range("D20").Value = '21.03'
$ws.Range("E20").Value = '  +7.32%  '

$ws.Range("D21").Value = '7.302'
$ws.Range("E21").Value = '  +0.45%  '

$ws.Range("D22").Value = '1.005'
$ws.Range("E22").Value = '  -0.48%  '

$ws.Range("D23").Value = '14.49'
$ws.Range("E23").Value = '  +2.44%  '

$ws.Range("D24").Value = '24.885.98'
$ws.Range("E24").Value = '  +1.37%  '

$ws.Range("D25").Value = '2.343'
$ws.Range("E25").Value = '  -0.29%  '

$ws.Range("D26").Value = '2.891'
$ws.Range("E26").Value = '  -3.60%  '

$ws.Range("B27").Value = 'HuobiToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D27").Value = '6.390'
$ws.Range("E27").Value = '  +21.93%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '23.10'
$ws.Range("E28").Value = '  +2.15%  '

$ws.Range("D29").Value = '162.30'
$ws.Range("E29").Value = '  +0.62%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '8.427'
$ws.Range("E30").Value = '  -0.23%  '

$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = '144.03'
$ws.Range("E31").Value = '  +5.46%  '

$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").Value = '2.281'
$ws.Range("E32").Value = '  +15.02%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.08786'
$ws.Range("E33").Value = '  +0.64%  '

$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D34").Value = '0.03191'
$ws.Range("E34").Value = '  +10.01%  '

$ws.Range("B35").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C35").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D35").Value = '1.824.00'
$ws.Range("E35").Value = '  -2.72%  '

$ws.Range("D36").Value = '7.208'
$ws.Range("E36").Value = '  -4.12%  '

$ws.Range("D37").Value = '1.032'
$ws.Range("E37").Value = '  +0.45%  '

$ws.Range("D38").Value = '0.2872'
$ws.Range("E38").Value = '  +6.02%  '

$ws.Range("D39").Value = '0.8445'
$ws.Range("E39").Value = '  +9.49%  '

$ws.Range("D40").Value = '10.87'
$ws.Range("E40").Value = '  +2.07%  '

$ws.Range("D41").Value = '0.09490'
$ws.Range("E41").Value = '  +4.35%  '

$ws.Range("D42").Value = '14.21'
$ws.Range("E42").Value = '  +1.27%  '

$ws.Range("D43").Value = '1.477'
$ws.Range("E43").Value = '  +1.65%  '

$ws.Range("D44").Value = '17.50'
$ws.Range("E44").Value = '  +5.32%  '

$ws.Range("D45").Value = '2.728'

$ws.Range("D46").Value = '0.7449'
$ws.Range("E46").Value = '  +4.69%  '

$ws.Range("E47").Value = '  +0.86%  '

$ws.Range("D48").Value = '1.383'
$ws.Range("E48").Value = '  +4.61%  '

$ws.Range("D49").Value = '1.001'
$ws.Range("E49").Value = '  -0.46%  '

$ws.Range("D50").Value = '140.90'

$ws.Range("D51").Value = '0.08399'
$ws.Range("E51").Value = '  +5.82%  '

# Restore default (unstyled) formatting on column D now that the text values are set,
# matching the workbook's original styling (no explicit NumberFormat on these cells).
$ws.Range("D2:D51").Style = "Normal"
